# LD50 task list - mark "define cricket character type" (row 21) and
# "relabel fireflies + firefly name refs to make more sense" (row 47) as
# done (Good/green), and flag the "multiple sets of dialogue for workers..."
# task (row 39) with a placeholder "AAAA" comment (Neutral/orange).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: define cricket character type -> done, assigned to Lex
$ws.Range("A21:D21").Style = "Good"
$ws.Range("B21").Value = "Lex"
$ws.Range("C21").Value = "done"

# Row 39: multiple sets of dialogue for workers... -> note "AAAA", assigned to Lex
$ws.Range("A39:D39").Style = "Neutral"
$ws.Range("B39").Value = "Lex"
$ws.Range("C39").Value = "AAAA"

# Row 47: relabel fireflies + firefly name refs to make more sense -> done
$ws.Range("A47:D47").Style = "Good"
$ws.Range("C47").Value = "done"

# Move the active selection to B14 (no scrolled "top left cell" pin)
$ws.Range("B14").Select() | Out-Null
